# Insert a new weekly record row at row 140 (shifts rows 140:184 down to 141:185)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(140).Insert()

$ws.Range("A140").Value = 8
$ws.Range("B140").Value = "Terminal La Palmera de La Serena"
$ws.Range("C140").Value = "Coquimbo"
$ws.Range("D140").Value = 44524
$ws.Range("E140").Value = 4
$ws.Range("F140").Value = 100112012
$ws.Range("G140").Value = "Espinaca"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 450
$ws.Range("L140").Value = 500
$ws.Range("M140").Value = 475
$ws.Range("N140").Value = "$/atado 300 a 500 gramos"
$ws.Range("O140").Value = "Provincia del Elquí"
$ws.Range("P140").Value = 950
$ws.Range("Q140").Value = 0.5
$ws.Range("R140").Value = "Hortaliza"
